# RotJ - up to final Level 2 boss, 312 ahead
#
# Fills in frame-count data for rows 45-50 on the "FrameCounts" sheet
# (Begin walljump -> Batman Appears, the run up to the final Level 2 boss),
# then moves the active selection to reflect where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FrameCounts")

# Row 45: Begin walljump
$ws.Range("C45").Value = 5492
$ws.Range("B45").Value = 5294
$ws.Range("A45").Value = "Begin walljump"

# Row 46: X = 188
$ws.Range("C46").Value = 5682
$ws.Range("B46").Value = 5475
$ws.Range("A46").Value = "X = 188"

# Row 47: X = 822
$ws.Range("C47").Value = 6210
$ws.Range("B47").Value = 5912
$ws.Range("A47").Value = "X = 822"

# Row 48: X = 1057
$ws.Range("C48").Value = 6265
$ws.Range("B48").Value = 5967
$ws.Range("A48").Value = "X = 1057"

# Row 49: Black screen
$ws.Range("C49").Value = 6474
$ws.Range("B49").Value = 6162
$ws.Range("A49").Value = "Black screen"

# Row 50: Batman Appears
$ws.Range("C50").Value = 6476
$ws.Range("B50").Value = 6164
$ws.Range("A50").Value = "Batman Appears"

# Leave the selection where the author left off editing.
$ws.Range("B51").Select()
